# Rebuild the per-row JSON-fragment formulas in column C of Sheet1 (FCC map
# scraper source). Two structural simplifications are applied throughout:
#   1. The leading array/object wrapper on row 2 drops the indentation
#      padding: `"[    {"` -> `" [{"`.
#   2. Every "subheader" node is flattened to a plain "name" node (the
#      `"subheader":{` wrapper is removed), and the matching closing-brace
#      cells that used to balance that wrapper (`] } } , {` -> `]  } , {`,
#      etc.) are tightened up to match the new, shallower nesting.
# Note: the workbook.xml `x15ac:absPath` (the author's local save-path
# breadcrumb) is not part of the Excel object model surfaced over COM
# (Workbook only exposes Name/Path/FullName/Worksheets/Sheets/ActiveSheet/
# Close) so it cannot be targeted from here; it is inert metadata, not
# worksheet content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Formula = '=" [{"'
$ws.Range("C72").Formula = '="""name"":"""&B72&""","'
$ws.Range("C104").Formula = '=""""&B104&"""]  } , {"'
$ws.Range("C105").Formula = '="""name"":"""&B105&""","'
$ws.Range("C110").Formula = '=""""&B110&"""]  } , {"'
$ws.Range("C111").Formula = '="""name"":"""&B111&""","'
$ws.Range("C130").Formula = '=""""&B130&"""]  } , {"'
$ws.Range("C131").Formula = '="""name"":"""&B131&""","'
$ws.Range("C135").Formula = '=""""&B135&"""]  } , {"'
$ws.Range("C136").Formula = '="""name"":"""&B136&""","'
$ws.Range("C241").Formula = '=""""&B241&"""]  } , {"'
$ws.Range("C242").Formula = '="""name"":"""&B242&""","'
$ws.Range("C256").Formula = '=""""&B256&"""]  } , {"'
$ws.Range("C257").Formula = '="""name"":"""&B257&""","'
$ws.Range("C275").Formula = '=""""&B275&"""]  } , {"'
$ws.Range("C276").Formula = '="""name"":"""&B276&""","'
$ws.Range("C284").Formula = '=""""&B284&"""]  } , {"'
$ws.Range("C285").Formula = '="""name"":"""&B285&""","'
$ws.Range("C290").Formula = '=""""&B290&"""]  } , {"'
$ws.Range("C291").Formula = '="""name"":"""&B291&""","'
$ws.Range("C313").Formula = '=""""&B313&"""]  } , {"'
$ws.Range("C314").Formula = '="""name"":"""&B314&""","'
$ws.Range("C324").Formula = '=""""&B324&"""]  } , {"'
$ws.Range("C325").Formula = '="""name"":"""&B325&""","'
$ws.Range("C330").Formula = '=""""&B330&"""]  } , {"'
$ws.Range("C331").Formula = '="""name"":"""&B331&""","'
$ws.Range("C333").Formula = '=""""&B333&"""] }  ]}  , {"'
$ws.Range("C337").Formula = '=""""&B337&"""]  } , {"'
$ws.Range("C338").Formula = '="""name"":"""&B338&""","'
$ws.Range("C340").Formula = '=""""&B340&"""]  } , {"'
$ws.Range("C341").Formula = '="""name"":"""&B341&""","'
$ws.Range("C347").Formula = '=""""&B347&"""]  } , {"'
$ws.Range("C348").Formula = '="""name"":"""&B348&""","'
$ws.Range("C350").Formula = '=""""&B350&"""]  } , {"'
$ws.Range("C351").Formula = '="""name"":"""&B351&""","'
$ws.Range("C357").Formula = '=""""&B357&"""]  } , {"'
$ws.Range("C358").Formula = '="""name"":"""&B358&""","'
$ws.Range("C360").Formula = '=""""&B360&"""] }  ]}  , {"'
$ws.Range("C365").Formula = '=""""&B365&"""]  } , {"'
$ws.Range("C366").Formula = '="""name"":"""&B366&""","'
$ws.Range("C368").Formula = '=""""&B368&"""]  } , {"'
$ws.Range("C369").Formula = '="""name"":"""&B369&""","'
$ws.Range("C375").Formula = '=""""&B375&"""]  } , {"'
$ws.Range("C376").Formula = '="""name"":"""&B376&""","'
$ws.Range("C378").Formula = '=""""&B378&"""]  } , {"'
$ws.Range("C379").Formula = '="""name"":"""&B379&""","'
$ws.Range("C386").Formula = '=""""&B386&"""]  } , {"'
$ws.Range("C387").Formula = '="""name"":"""&B387&""","'
$ws.Range("C394").Formula = '=""""&B394&"""]  } , {"'
$ws.Range("C395").Formula = '="""name"":"""&B395&""","'
$ws.Range("C397").Formula = '=""""&B397&"""] }  ]}  , {"'
$ws.Range("C405").Formula = '=""""&B405&"""] }  ]}  , {"'
$ws.Range("C411").Formula = '=""""&B411&"""]  } , {"'
$ws.Range("C412").Formula = '="""name"":"""&B412&""","'
